# Add the new AVTA Horticulture / Agribusiness promotion course rows
# (rows 2-8) to the "courses" sheet, matching the offshore Region3
# promotions workbook update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{
        R = 2
        A = "AHC30716"
        B = "110597F"
        C = "HORTICULTURE"
        D = "CERTIFICATE III IN HORTICULTURE"
        E = 52
        H = "44 wks Tuition + 8 wks Break"
        I = 10200
        J = "10,000 tuition fee + 200 handling fee"
        M = "TAS"
        RR = "PROMOTIONS VALID UNTIL 28TH FEBRUARY 2023 "
        Height = 45
        WrapABD = $false
    },
    @{
        R = 3
        A = "AHC40416"
        B = "110598E"
        C = "HORTICULTURE"
        D = "CERTIFICATE IV IN HORTICULTURE"
        E = 52
        H = "44 wks Tuition + 8 wks Break"
        I = 11200
        J = "11,000 tuition fee + 200 handling fee"
        M = "TAS"
        RR = "PROMOTIONS VALID UNTIL 28TH FEBRUARY 2023 "
        Height = 45
        WrapABD = $false
    },
    @{
        R = 4
        A = "AHC51422"
        B = "110774E"
        C = "MANAGEMENT"
        D = "DIPLOMA OF AGRIBUSINESS MANAGEMENT"
        E = 52
        H = "44 wks Tuition + 8 wks Break"
        I = 12200
        J = "12,000 tuition fee + 200 handling fee"
        M = "TAS"
        RR = "PROMOTIONS VALID UNTIL 28TH FEBRUARY 2023 "
        Height = 45
        WrapABD = $false
    },
    @{
        R = 5
        A = "AHC30716 / AHC40416"
        B = "110597F / 110598E"
        C = "PACKAGES"
        D = "CERTIFICATE III IN HORTICULTURE +`nCERTIFICATE IV IN HORTICULTURE"
        E = 104
        H = "88 wks Tuition + 16 wks Break"
        I = 22200
        J = "22,000 tuition fee + 200 handling fee"
        M = "TAS"
        RR = "PROMOTIONS VALID UNTIL 28TH FEBRUARY 2023 "
        Height = 45
        WrapABD = $true
    },
    @{
        R = 6
        A = "AHC30716 / AHC51422"
        B = "110597F / 110774E"
        C = "PACKAGES"
        D = "CERTIFICATE III IN HORTICULTURE +`nDIPLOMA OF AGRIBUSINESS MANAGEMENT"
        E = 104
        H = "88 wks Tuition + 16 wks Break"
        I = 21200
        J = "21,000 tuition fee + 200 handling fee"
        M = "TAS"
        RR = "PROMOTIONS VALID UNTIL 28TH FEBRUARY 2023 "
        Height = 45
        WrapABD = $true
    },
    @{
        R = 7
        A = "AHC40416 / AHC51422"
        B = "110598E / 110774E"
        C = "PACKAGES"
        D = "CERTIFICATE IV IN HORTICULTURE +`nDIPLOMA OF AGRIBUSINESS MANAGEMENT"
        E = 104
        H = "88 wks Tuition + 16 wks Break"
        I = 22200
        J = "22,000 tuition fee + 200 handling fee"
        M = "TAS"
        RR = "PROMOTIONS VALID UNTIL 28TH FEBRUARY 2023 "
        Height = 45
        WrapABD = $true
    },
    @{
        R = 8
        A = "AHC30716 / AHC40416 / AHC51422"
        B = "110597F / 110598E / 110774E"
        C = "PACKAGES"
        D = "CERTIFICATE III IN HORTICULTURE +`nCERTIFICATE IV IN HORTICULTURE +`nDIPLOMA OF AGRIBUSINESS MANAGEMENT"
        E = 156
        H = "132 wks Tuition + 24 wks Break"
        I = 30200
        J = "30,000 tuition fee + 200 handling fee"
        M = "TAS"
        RR = "PROMOTIONS VALID UNTIL 28TH FEBRUARY 2023 "
        Height = 75
        WrapABD = $true
    }
)

foreach ($row in $rows) {
    $r = $row.R

    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E

    $hCell = $ws.Cells.Item($r, 8)
    $hCell.Value = $row.H
    $hCell.WrapText = $true

    $iCell = $ws.Cells.Item($r, 9)
    $iCell.Value = $row.I

    $jCell = $ws.Cells.Item($r, 10)
    $jCell.Value = $row.J

    $ws.Range($iCell, $jCell).NumberFormat = "#,##0"
    $jCell.WrapText = $true

    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 18).Value = $row.RR

    if ($row.WrapABD) {
        $ws.Cells.Item($r, 1).WrapText = $true
        $ws.Cells.Item($r, 2).WrapText = $true
        $ws.Cells.Item($r, 4).WrapText = $true
    }

    $ws.Rows($r).RowHeight = $row.Height
}

$ws.Range("F17").Select() | Out-Null
